$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.274.97"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.18"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.13"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4675"
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06553"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.83"
$ws.Range("E10").Value = "  +6.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07880"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.20"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.868.78"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.150"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6760"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.09"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.280.05"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.514"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.68"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.111.85"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007276"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.182"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.311"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.28"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.13"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.910"
$ws.Range("E28").Value = "  -5.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.353"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09678"
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.429"
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04704"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.106"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7044"
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("E39").Value = "  -4.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.533"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.32"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.942"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8479"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4186"
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.16"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9993"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.212"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.281"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "942.61"
$ws.Range("E49").Value = "  -4.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.14"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("E51").Value = "  -4.30%  "
